# Applies the "Updating github with sharelatex" change to the
# cc-aggregate sheet:
#   - adds a bold "FULL SIZE" label at E15
#   - adds a new F12-H12 computation at E20
#   - adds a new "STAGE 1 / CUM" cumulative-sum table in A26:B36
#   - moves the active selection to B37

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E15: new bold "FULL SIZE" header cell ------------------------------
$ws.Range("E15").Value = "FULL SIZE"
$ws.Range("E15").Font.Bold = $true

# --- E20: new formula cell ----------------------------------------------
$ws.Range("E20").Formula = "=F12-H12"

# --- A26:B26: new section header ("STAGE 1" / "CUM") --------------------
$ws.Range("A26").Value = "STAGE 1"
$ws.Range("B26").Value = "CUM"

# --- A27:B36: new cumulative-sum table -----------------------------------
# Column A inherits a text (@) number format from the sheet's column
# style, so reset A27:B36 back to the default "Normal" style first;
# otherwise the numbers we enter would be stored as text.
$ws.Range("A27:B36").Style = "Normal"

$ws.Range("A27").Value = 116649
$ws.Range("B27").Formula = "=A27"

$ws.Range("A28").Value = 145142
$ws.Range("B28").Formula = "=A27+A28"

$ws.Range("A29").Value = 179943
$ws.Range("B29").Formula = "=B28+A29"

$ws.Range("A30").Value = 229916
$ws.Range("A31").Value = 283874
$ws.Range("A32").Value = 236082
$ws.Range("A33").Value = 238769
$ws.Range("A34").Value = 351397
$ws.Range("A35").Value = 182947
$ws.Range("A36").Value = 162565

# Assign the relative formula to the whole B30:B36 block in one call so
# the engine fills it down as a single shared formula (matching how the
# workbook was originally authored), rather than as per-cell formulas.
$ws.Range("B30:B36").Formula = "=B29+A30"

# --- move the active cell/selection to B37 (below the new table) --------
$ws.Range("B37").Select()
